$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.589.63"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "2.352.63"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.634"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.17%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0924"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.46"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.993"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.23%  "
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.04%  "
$ws.Range("D16").Value = "2.709.60"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").Value = "2.366.43"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").Value = "42.563.45"
$ws.Range("E18").Value = "  -1.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.76%  "
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "256.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.08%  "
$ws.Range("E24").Value = "  -3.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.87%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("E29").Value = "  +3.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0888"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -10.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.125"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +18.57%  "
$ws.Range("E36").Value = "  -1.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0362"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.237"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.32%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "112.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -12.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.15%  "
